$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Homework")
$c = $ws.Cells.Item(4,2)
Write-Output ("Font.Size=" + $c.Font.Size)
Write-Output ("Font.Name=" + $c.Font.Name)
Write-Output ("Interior.ColorIndex=" + $c.Interior.ColorIndex)
Write-Output ("Interior.Color=" + $c.Interior.Color)
Write-Output ("StyleID?")
